$wb = $excel.ActiveWorkbook

# --- 1. Status text update: "Ready for handoff" -> "In Translation" ---
# Overview sheet keeps per-locale status in columns E (zh-cn) and F (de-de).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"

# The zh-cn / de-de detail sheets keep the same status in column C.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"

# --- 2. Narrow the status columns ---
# ColumnWidth is quantized by the host to a 1/6-character pixel grid, so
# 12.5 is the representable width closest to the recorded target (~13.41).
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
